$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, $cellRef, $value)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws = $wb.Worksheets.Item("Commitments")
Set-TextValue $ws "B9" "CMT - 000950"
Set-TextValue $ws "H9" "6/13/2022"

$ws = $wb.Worksheets.Item("CapitalCall")
Set-TextValue $ws "B2" "CC-0662"
Set-TextValue $ws "C2" "DD-0224"
Set-TextValue $ws "I2" "6/17/2022"
Set-TextValue $ws "J2" "6/18/2022"
Set-TextValue $ws "B3" "CC-0663"
Set-TextValue $ws "C3" "DD-0224"
Set-TextValue $ws "I3" "6/17/2022"
Set-TextValue $ws "J3" "6/18/2022"
Set-TextValue $ws "B4" "CC-0664"
Set-TextValue $ws "C4" "DD-0224"
Set-TextValue $ws "I4" "6/17/2022"
Set-TextValue $ws "J4" "6/18/2022"
Set-TextValue $ws "L4" "6/17/2022"
Set-TextValue $ws "B5" "CC-0665"
Set-TextValue $ws "C5" "DD-0224"
Set-TextValue $ws "I5" "6/17/2022"
Set-TextValue $ws "J5" "6/18/2022"
Set-TextValue $ws "B6" "CC-0666"
Set-TextValue $ws "C6" "DD-0224"
Set-TextValue $ws "I6" "6/17/2022"
Set-TextValue $ws "J6" "6/18/2022"

$ws = $wb.Worksheets.Item("FundDrawdown")
Set-TextValue $ws "B2" "DD-0224"
Set-TextValue $ws "E2" "6/17/2022"
Set-TextValue $ws "F2" "6/18/2022"

$ws = $wb.Worksheets.Item("FundDistribution")
Set-TextValue $ws "B2" "FD-0233"

$ws = $wb.Worksheets.Item("InvestorDistribution")
Set-TextValue $ws "B2" "FD-0233"
Set-TextValue $ws "C2" "ID-0648"
Set-TextValue $ws "D2" "CMT - 000952"
Set-TextValue $ws "E2" "15000000"
Set-TextValue $ws "F2" "3000000"
Set-TextValue $ws "G2" "200000.000"
Set-TextValue $ws "H2" "200000.000"
Set-TextValue $ws "I2" "200000.000"
Set-TextValue $ws "J2" "3600000.000"
Set-TextValue $ws "B3" "FD-0233"
Set-TextValue $ws "C3" "ID-0649"
Set-TextValue $ws "D3" "CMT - 000951"
Set-TextValue $ws "E3" "10000000"
Set-TextValue $ws "F3" "2000000"
Set-TextValue $ws "G3" "133333.333"
Set-TextValue $ws "H3" "133333.333"
Set-TextValue $ws "I3" "133333.333"
Set-TextValue $ws "J3" "2399999.999"
Set-TextValue $ws "B4" "FD-0233"
Set-TextValue $ws "C4" "ID-0650"
Set-TextValue $ws "D4" "CMT - 000950"
Set-TextValue $ws "B5" "FD-0233"
Set-TextValue $ws "C5" "ID-0651"
Set-TextValue $ws "D5" "CMT - 000953"
Set-TextValue $ws "E5" "25000000"
Set-TextValue $ws "F5" "5000000"
Set-TextValue $ws "G5" "333333.333"
Set-TextValue $ws "H5" "333333.333"
Set-TextValue $ws "I5" "333333.333"
Set-TextValue $ws "J5" "5999999.999"
Set-TextValue $ws "B6" "FD-0233"
Set-TextValue $ws "C6" "ID-0652"
Set-TextValue $ws "D6" "CMT - 000954"
Set-TextValue $ws "E6" "5000000"
Set-TextValue $ws "F6" "1000000"
Set-TextValue $ws "G6" "66666.667"
Set-TextValue $ws "H6" "66666.667"
Set-TextValue $ws "I6" "66666.667"
Set-TextValue $ws "J6" "1200000.001"
